# "Attempting to implement auto highlighting"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "test" row (row 2) content ---
# B2: ausgespielter Flow/FAQ sample text
$ws.Range("B2").Value = "Mediterranean Tuna Steaks"

# D2: sample bot answer text (multi-line)
$ws.Range("D2").Value = "`n`nWhat about the ups?`nThat's good.`nHaha, well, let's just say you don't want to be involved with i"

# E2: feedback on how well the flow worked
$ws.Range("E2").Value = "Es hat nicht gut geklappt, da der Bot das Anliegen nicht richtig erkannt hat. FAQ xy wäre besser.0420000000"

# F2 (new column): a "Kommentar" cell, highlighted in yellow to test auto highlighting
$ws.Range("F2").Value = "Details"
$ws.Range("F2").Interior.Color = 65535

# Row 2 grew taller to fit the new content
$ws.Rows.Item(2).RowHeight = 144

# --- Row 7 gets a stray numeric test value in column F ---
$ws.Range("F7").Value = 420000000

# --- Window / selection bookkeeping ---
$ws.Range("A5").Select()
$wb.Windows.Item(1).WindowState = -4140
